$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "is_active" column (E2:E13) currently holds =TRUE() formulas (numeric
# boolean results). Replace them with the literal text string "TRUE" while
# keeping the existing cell formatting (style) intact.
#
# Simply assigning the string "TRUE" via .Value/.Formula gets auto-coerced
# by Excel into a real Boolean cell, so instead we build the text value in
# a scratch cell (as a formula result, so it is a genuine string) and copy
# it across with Paste Special “Values”, which carries over the text value
# only and leaves each destination cell's existing number format/style
# untouched.
$helper = $ws.Cells.Item(20, 20)
$helper.Formula = '="TRUE"'
$helper.Copy()

for ($r = 2; $r -le 13; $r++) {
    $target = $ws.Cells.Item($r, 5)
    $target.PasteSpecial(-4163)  # xlPasteValues
}

$helper.Clear()

# Match the saved selection state: E2 active, E2:E13 selected.
$ws.Range("E2:E13").Select()
